$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns are treated as plain text so that
# values such as "1.001" or "27.482.02" are not auto-converted into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.482.02'
$ws.Range("E2").Value = '  -0.52%  '

$ws.Range("D3").Value = '1.824.23'
$ws.Range("E3").Value = '  -0.95%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '311.82'
$ws.Range("E5").Value = '  -0.11%  '

$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").Value = '0.4235'
$ws.Range("E7").Value = '  -0.46%  '

$ws.Range("D8").Value = '0.3625'
$ws.Range("E8").Value = '  +0.46%  '

$ws.Range("D9").Value = '0.07174'
$ws.Range("E9").Value = '  -1.81%  '

$ws.Range("D10").Value = '0.8580'
$ws.Range("E10").Value = '  -1.90%  '

$ws.Range("D11").Value = '20.57'
$ws.Range("E11").Value = '  +0.01%  '

$ws.Range("D12").Value = '1.882.71'
$ws.Range("E12").Value = '  +3.59%  '

$ws.Range("D13").Value = '5.400'
$ws.Range("E13").Value = '  +1.48%  '

$ws.Range("D14").Value = '6.470'
$ws.Range("E14").Value = '  -0.26%  '

$ws.Range("D15").Value = '0.06914'
$ws.Range("E15").Value = '  -0.78%  '

$ws.Range("D16").Value = '1.004'
$ws.Range("E16").Value = '  +0.04%  '

$ws.Range("D17").Value = '80.14'
$ws.Range("E17").Value = '  +1.13%  '

$ws.Range("D18").Value = '0.000008865'
$ws.Range("E18").Value = '  -0.69%  '

$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  -0.09%  '

$ws.Range("D20").Value = '15.30'
$ws.Range("E20").Value = '  -0.06%  '

$ws.Range("D21").Value = '27.272.04'
$ws.Range("E21").Value = '  -1.65%  '

$ws.Range("D22").Value = '5.115'
$ws.Range("E22").Value = '  +3.10%  '

$ws.Range("D23").Value = '10.92'
$ws.Range("E23").Value = '  +5.97%  '

$ws.Range("D24").Value = '2.027.41'
$ws.Range("E24").Value = '  -0.73%  '

$ws.Range("D25").Value = '1.985'
$ws.Range("E25").Value = '  -0.22%  '

$ws.Range("D26").Value = '154.99'
$ws.Range("E26").Value = '  -0.33%  '

$ws.Range("E27").Value = '  +1.14%  '

$ws.Range("D28").Value = '5.131'
$ws.Range("E28").Value = '  -1.25%  '

$ws.Range("D29").Value = '113.79'
$ws.Range("E29").Value = '  -4.61%  '

$ws.Range("D30").Value = '1.802'
$ws.Range("E30").Value = '  -3.77%  '

$ws.Range("D31").Value = '0.08839'
$ws.Range("E31").Value = '  -0.22%  '

$ws.Range("D32").Value = '2.985'

$ws.Range("D33").Value = '0.7425'
$ws.Range("E33").Value = '  -2.06%  '

$ws.Range("E34").Value = '  +0.82%  '

$ws.Range("E35").Value = '  -0.32%  '

$ws.Range("E36").Value = '  +0.10%  '

$ws.Range("E37").Value = '  -1.39%  '

$ws.Range("D38").Value = '0.05278'
$ws.Range("E38").Value = '  -2.56%  '

$ws.Range("D39").Value = '0.01924'
$ws.Range("E39").Value = '  +0.08%  '

$ws.Range("E40").Value = '  -1.16%  '

$ws.Range("D41").Value = '0.5038'
$ws.Range("E41").Value = '  -0.23%  '

$ws.Range("E42").Value = '  -0.92%  '

$ws.Range("D43").Value = '6.478'
$ws.Range("E43").Value = '  -0.80%  '

$ws.Range("D44").Value = '8.283'
$ws.Range("E44").Value = '  -0.92%  '

$ws.Range("D45").Value = '10.38'
$ws.Range("E45").Value = '  +0.09%  '

$ws.Range("D46").Value = '105.41'
$ws.Range("E46").Value = '  -0.53%  '

$ws.Range("D47").Value = '0.06449'
$ws.Range("E47").Value = '  -1.50%  '

$ws.Range("D48").Value = '0.4663'
$ws.Range("E48").Value = '  +0.93%  '

$ws.Range("E49").Value = '  +0.02%  '

$ws.Range("D50").Value = '1.612'
$ws.Range("E50").Value = '  -1.32%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '63.23'
$ws.Range("E51").Value = '  -1.54%  '
